$wb = $excel.ActiveWorkbook

# --- Test Cases sheet: disable (Runmode = N) the ClaimMojioCancel and
#     InvalidIMEISubmit test cases ---
$wsCases = $wb.Worksheets.Item("Test Cases")
$wsCases.Range("C5").Value = "N"
$wsCases.Range("C6").Value = "N"

# --- Test Steps sheet: clear the recorded Result1 for the steps that
#     belong to the now-disabled test cases (ClaimMojioCancel rows 20-25,
#     InvalidIMEISubmit rows 26-31) ---
$wsSteps = $wb.Worksheets.Item("Test Steps")
$wsSteps.Range("H20:H31").ClearContents()

# refresh the still-enabled ClaimMojioSubmit rows' Result1 values
for ($r = 14; $r -le 19; $r++) {
    $wsSteps.Cells.Item($r, 8).Value = "PASS"
}

# move the view/selection on the Test Steps sheet back to the left, onto
# the Result1 header cell
$wsSteps.Range("H1").Select()

# --- Restore tab selection to "Test Cases" (it was left on
#     "InvalidIMEISubmit" before), placing the cursor on D5 ---
$wsCases.Activate()
$wsCases.Range("D5").Select()
